$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.25756710338925
$ws.Range("C2").Value = 5.342309548077227
$ws.Range("D2").Value = 9.297245342656394
$ws.Range("E2").Value = 13.7335740058896
$ws.Range("F2").Value = 32.52346818206428
$ws.Range("I2").Value = 22.24974157687664
$ws.Range("J2").Value = 9.921764148115141
$ws.Range("K2").Value = 10.47487306526676
$ws.Range("N2").Value = 18.85652582045357
$ws.Range("O2").Value = 24.48923716574427
$ws.Range("B3").Value = 9.984085959012322
$ws.Range("C3").Value = 5.147060776718794
$ws.Range("D3").Value = 9.234468675092797
$ws.Range("E3").Value = 13.66831861141287
$ws.Range("F3").Value = 32.56530234486023
$ws.Range("I3").Value = 22.33373777365707
$ws.Range("J3").Value = 9.927571469992628
$ws.Range("K3").Value = 10.29322046947333
$ws.Range("N3").Value = 18.91345717709833
$ws.Range("O3").Value = 24.55969453908656
$ws.Range("B4").Value = 9.813959927715334
$ws.Range("C4").Value = 5.024042825955796
$ws.Range("D4").Value = 9.197417890228325
$ws.Range("E4").Value = 13.63099248755584
$ws.Range("F4").Value = 32.59838166921287
$ws.Range("I4").Value = 22.38922204329289
$ws.Range("J4").Value = 9.9326860440286
$ws.Range("K4").Value = 10.18169473288701
$ws.Range("N4").Value = 18.95006846841519
$ws.Range("O4").Value = 24.60776970458086
$ws.Range("B5").Value = 9.744183429228711
$ws.Range("C5").Value = 4.973201039009396
$ws.Range("D5").Value = 9.182708343453371
$ws.Range("E5").Value = 13.61648340286945
$ws.Range("F5").Value = 32.61371867737716
$ws.Range("I5").Value = 22.4128151439122
$ws.Range("J5").Value = 9.935160187147012
$ws.Range("K5").Value = 10.13630764894352
$ws.Range("N5").Value = 18.96540528594443
$ws.Range("O5").Value = 24.62856941878221
$ws.Range("B6").Value = 9.732573209732053
$ws.Range("C6").Value = 4.964718264781778
$ws.Range("D6").Value = 9.180289699009297
$ws.Range("E6").Value = 13.61411689918657
$ws.Range("F6").Value = 32.61637745453412
$ws.Range("I6").Value = 22.416792106563
$ws.Range("J6").Value = 9.935594577258909
$ws.Range("K6").Value = 10.12877652901309
$ws.Range("N6").Value = 18.96797719933506
$ws.Range("O6").Value = 24.63209614194823
$ws.Range("B7").Value = 9.813020569770535
$ws.Range("C7").Value = 5.023359926888494
$ws.Range("D7").Value = 9.197217920337271
$ws.Range("E7").Value = 13.63079395650552
$ws.Range("F7").Value = 32.59858099417683
$ws.Range("I7").Value = 22.3895362493383
$ws.Range("J7").Value = 9.932717831942432
$ws.Range("K7").Value = 10.1810823042552
$ws.Range("N7").Value = 18.9502736143913
$ws.Range("O7").Value = 24.60804532535947
$ws.Range("B8").Value = 10.16379096040405
$ws.Range("C8").Value = 5.275682272926438
$ws.Range("D8").Value = 9.27529789075165
$ws.Range("E8").Value = 13.7105127720389
$ws.Range("F8").Value = 32.5363571082076
$ws.Range("I8").Value = 22.27789161465178
$ws.Range("J8").Value = 9.923445341818113
$ws.Range("K8").Value = 10.41227184375979
$ws.Range("N8").Value = 18.87581294114672
$ws.Range("O8").Value = 24.51253063482084
$ws.Range("B9").Value = 10.82976924833981
$ws.Range("C9").Value = 5.742607939962167
$ws.Range("D9").Value = 9.439662774317982
$ws.Range("E9").Value = 13.88800279789963
$ws.Range("F9").Value = 32.47306548654721
$ws.Range("I9").Value = 22.09000344145979
$ws.Range("J9").Value = 9.917528784546327
$ws.Range("K9").Value = 10.86320818063559
$ws.Range("N9").Value = 18.74287334184322
$ws.Range("O9").Value = 24.36350302943195
$ws.Range("B10").Value = 11.30033681783912
$ws.Range("C10").Value = 6.065149207830316
$ws.Range("D10").Value = 9.566432911740383
$ws.Range("E10").Value = 14.03046189438903
$ws.Range("F10").Value = 32.46242533222507
$ws.Range("I10").Value = 21.97091672658866
$ws.Range("J10").Value = 9.920624540936094
$ws.Range("K10").Value = 11.18978166969497
$ws.Range("N10").Value = 18.65309574235725
$ws.Range("O10").Value = 24.27745002746461
$ws.Range("B11").Value = 11.50934082294158
$ws.Range("C11").Value = 6.206822168571382
$ws.Range("D11").Value = 9.625220308046508
$ws.Range("E11").Value = 14.09769235535224
$ws.Range("F11").Value = 32.46536903938006
$ws.Range("I11").Value = 21.92086331595026
$ws.Range("J11").Value = 9.923639104852029
$ws.Range("K11").Value = 11.3366700977145
$ws.Range("N11").Value = 18.61395060943248
$ws.Range("O11").Value = 24.24341141619806
$ws.Range("B12").Value = 11.58768337372807
$ws.Range("C12").Value = 6.259701194822827
$ws.Range("D12").Value = 9.64762613965555
$ws.Range("E12").Value = 14.1234819591572
$ws.Range("F12").Value = 32.46760147230432
$ws.Range("I12").Value = 21.90250244880735
$ws.Range("J12").Value = 9.925010559799183
$ws.Range("K12").Value = 11.39200130176248
$ws.Range("N12").Value = 18.59936988239853
$ws.Range("O12").Value = 24.23125774153401
$ws.Range("B13").Value = 11.57084773038119
$ws.Range("C13").Value = 6.248347595968697
$ws.Range("D13").Value = 9.642794490392495
$ws.Range("E13").Value = 14.11791328031144
$ws.Range("F13").Value = 32.46707100132713
$ws.Range("I13").Value = 21.9064303902375
$ws.Range("J13").Value = 9.924704985830674
$ws.Range("K13").Value = 11.3800985339567
$ws.Range("N13").Value = 18.60249932634908
$ws.Range("O13").Value = 24.2338424949918
$ws.Range("B14").Value = 11.51580260617692
$ws.Range("C14").Value = 6.211188204740031
$ws.Range("D14").Value = 9.627060852604966
$ws.Range("E14").Value = 14.0998075543064
$ws.Range("F14").Value = 32.46553031434637
$ws.Range("I14").Value = 21.91934086130624
$ws.Range("J14").Value = 9.923747334250184
$ws.Range("K14").Value = 11.34122837450082
$ws.Range("N14").Value = 18.61274618658626
$ws.Range("O14").Value = 24.24239676681063
$ws.Range("B15").Value = 11.4819792236837
$ws.Range("C15").Value = 6.188325621815856
$ws.Range("D15").Value = 9.617441844245386
$ws.Range("E15").Value = 14.08875982641279
$ws.Range("F15").Value = 32.46473209726103
$ws.Range("I15").Value = 21.92732618506333
$ws.Range("J15").Value = 9.923190651327063
$ws.Range("K15").Value = 11.31737972285775
$ws.Range("N15").Value = 18.61905425743197
$ws.Range("O15").Value = 24.24773239984521
$ws.Range("B16").Value = 11.28656952229433
$ws.Range("C16").Value = 6.055785080066247
$ws.Range("D16").Value = 9.562612120854729
$ws.Range("E16").Value = 14.02611551331223
$ws.Range("F16").Value = 32.46238946777718
$ws.Range("I16").Value = 21.97427080750683
$ws.Range("J16").Value = 9.92045976547911
$ws.Range("K16").Value = 11.18014403618522
$ws.Range("N16").Value = 18.65568798740458
$ws.Range("O16").Value = 24.27977744166406
$ws.Range("B17").Value = 11.16534047827749
$ws.Range("C17").Value = 5.973149328450765
$ws.Range("D17").Value = 9.529250730623438
$ws.Range("E17").Value = 13.9882938794545
$ws.Range("F17").Value = 32.46294510629749
$ws.Range("I17").Value = 22.00412558333334
$ws.Range("J17").Value = 9.919195137530497
$ws.Range("K17").Value = 11.09548854917281
$ws.Range("N17").Value = 18.67859499361964
$ws.Range("O17").Value = 24.30074531614018
$ws.Range("B18").Value = 11.09514056666591
$ws.Range("C18").Value = 5.925146035354611
$ws.Range("D18").Value = 9.510168477661464
$ws.Range("E18").Value = 13.9667692457728
$ws.Range("F18").Value = 32.46399750485921
$ws.Range("I18").Value = 22.02168497717158
$ws.Range("J18").Value = 9.918619020833102
$ws.Range("K18").Value = 11.04664257124748
$ws.Range("N18").Value = 18.69193011652229
$ws.Range("O18").Value = 24.31328612328837
$ws.Range("B19").Value = 11.07129338897513
$ws.Range("C19").Value = 5.908813038213169
$ws.Range("D19").Value = 9.50372632029779
$ws.Range("E19").Value = 13.95952130302039
$ws.Range("F19").Value = 32.46447972588722
$ws.Range("I19").Value = 22.02769685089961
$ws.Range("J19").Value = 9.918449968118106
$ws.Range("K19").Value = 11.03007927034585
$ws.Range("N19").Value = 18.69647260764174
$ws.Range("O19").Value = 24.31761472510021
$ws.Range("B20").Value = 11.17829495107112
$ws.Range("C20").Value = 5.981995368168603
$ws.Range("D20").Value = 9.532791221736499
$ws.Range("E20").Value = 13.99229644542121
$ws.Range("F20").Value = 32.46281012368404
$ws.Range("I20").Value = 22.00090735808561
$ws.Range("J20").Value = 9.919314111504599
$ws.Range("K20").Value = 11.10451663217206
$ws.Range("N20").Value = 18.67613999157157
$ws.Range("O20").Value = 24.29846349372521
$ws.Range("B21").Value = 11.53199304486845
$ws.Range("C21").Value = 6.222124015631667
$ws.Range("D21").Value = 9.63167841975792
$ws.Range("E21").Value = 14.10511680972767
$ws.Range("F21").Value = 32.4659525337561
$ws.Range("I21").Value = 21.91553263481165
$ws.Range("J21").Value = 9.924022389380456
$ws.Range("K21").Value = 11.35265381509544
$ws.Range("N21").Value = 18.60972985744825
$ws.Range("O21").Value = 24.23986418286226
$ws.Range("B22").Value = 11.75844612987041
$ws.Range("C22").Value = 6.374560937936648
$ws.Range("D22").Value = 9.697140319404793
$ws.Range("E22").Value = 14.18077144616177
$ws.Range("F22").Value = 32.47451987057931
$ws.Range("I22").Value = 21.86319377632717
$ws.Range("J22").Value = 9.928438995859496
$ws.Range("K22").Value = 11.51310009458761
$ws.Range("N22").Value = 18.56774106048747
$ws.Range("O22").Value = 24.2058567992193
$ws.Range("B23").Value = 11.63803783299683
$ws.Range("C23").Value = 6.29362708750129
$ws.Range("D23").Value = 9.662131341560791
$ws.Range("E23").Value = 14.14022349853714
$ws.Range("F23").Value = 32.46935206470013
$ws.Range("I23").Value = 21.89081127323377
$ws.Range("J23").Value = 9.9259595959218
$ws.Range("K23").Value = 11.42764127074017
$ws.Range("N23").Value = 18.59002223679715
$ws.Range("O23").Value = 24.22361409963854
$ws.Range("B24").Value = 11.17243979710251
$ws.Range("C24").Value = 5.977997611796996
$ws.Range("D24").Value = 9.531190260269376
$ws.Range("E24").Value = 13.99048620020822
$ws.Range("F24").Value = 32.46286886610272
$ws.Range("I24").Value = 22.00236108586139
$ws.Range("J24").Value = 9.91925985314692
$ws.Range("K24").Value = 11.1004355813427
$ws.Range("N24").Value = 18.6772493821337
$ws.Range("O24").Value = 24.29949359158744
$ws.Range("B25").Value = 10.65252959041889
$ws.Range("C25").Value = 5.619666703328252
$ws.Range("D25").Value = 9.394078979813166
$ws.Range("E25").Value = 13.83780537589684
$ws.Range("F25").Value = 32.48389137838569
$ws.Range("I25").Value = 22.13750560257462
$ws.Range("J25").Value = 9.917819324322402
$ws.Range("K25").Value = 10.74182011742925
$ws.Range("N25").Value = 18.77744522385475
$ws.Range("O25").Value = 24.39971006594965
